# Map032 scene update: insert a new (blank) column C, shifting the
# existing column C data (header + entries) one column to the right
# into column D. Column B is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()
